$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the contents of A3 and C3 while keeping their styles
$ws.Range("A3").ClearContents()
$ws.Range("C3").ClearContents()

# Update the selection to match what the user had selected (A3:C3) when
# the contents were cleared, with C3 as the active cell.
$ws.Range("A3:C3").Select()
